$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 5 placeholder rows (8-12) previously had no data, just a blank
# custom row format. Clear that row-level formatting first so we can
# populate them with real validation-test data, each row missing (or
# wrong in) exactly one field, highlighted in red.
$ws.Rows("8:12").ClearFormats()

# Row 8: ServiceNow ID missing
$ws.Range("H8").Value = 1
$ws.Range("J8").Value = "09876543210987654322abcdefABCDEF"
$ws.Range("L8").Value = 1
$ws.Range("C8").Interior.Color = 255

# Row 9: Participant Count missing
$ws.Range("C9").Value = "09876543210987654323abcdefABCDEF"
$ws.Range("J9").Value = "09876543210987654323abcdefABCDEF"
$ws.Range("L9").Value = 1
$ws.Range("H9").Interior.Color = 255

# Row 10: Title missing
$ws.Range("C10").Value = "09876543210987654324abcdefABCDEF"
$ws.Range("H10").Value = 1
$ws.Range("L10").Value = 1
$ws.Range("J10").Interior.Color = 255

# Row 11: Participant Count result missing
$ws.Range("C11").Value = "09876543210987654325abcdefABCDEF"
$ws.Range("H11").Value = 1
$ws.Range("J11").Value = "09876543210987654325abcdefABCDEF"
$ws.Range("L11").Interior.Color = 255

# Row 12: Participant Count result has a wrong value
$ws.Range("C12").Value = "09876543210987654326abcdefABCDEF"
$ws.Range("H12").Value = 1
$ws.Range("J12").Value = "09876543210987654326abcdefABCDEF"
$ws.Range("L12").Value = "WRONG VALUE"
$ws.Range("L12").Interior.Color = 255

# Move the active selection to L11
$ws.Range("L11").Select()
